$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 410270.53
$ws.Range("I70").Value = 929066.6
$ws.Range("J70").Value = 2645
$ws.Range("K70").Value = 2787199.8
$ws.Range("L70").Value = 7935
$ws.Range("M70").Value = -2786929.8
$ws.Range("N70").Value = -8475
$ws.Range("H73").Value = 410270.53
$ws.Range("I73").Value = 929066.6
$ws.Range("J73").Value = 2645
$ws.Range("K73").Value = 2787199.8
$ws.Range("L73").Value = 7935
$ws.Range("M73").Value = -2786263.8
$ws.Range("N73").Value = -9807
$ws.Range("H80").Value = 1137571.2
$ws.Range("I80").Value = 1895628.4
$ws.Range("K80").Value = 5686885.199999999
$ws.Range("M80").Value = -5685887.199999999
$ws.Range("H83").Value = 1137571.2
$ws.Range("I83").Value = 1895628.4
$ws.Range("K83").Value = 17060655.6
$ws.Range("M83").Value = -17055663.6
$ws.Range("H111").Value = 76015.5
$ws.Range("I111").Value = 90000
$ws.Range("J111").Value = 71354
$ws.Range("K111").Value = 270000
$ws.Range("L111").Value = 214062
$ws.Range("M111").Value = -266933
$ws.Range("N111").Value = -220196
$ws.Range("H132").Value = 3605.8428
$ws.Range("I132").Value = 3468.8823
$ws.Range("K132").Value = 10406.6469
$ws.Range("M132").Value = -7876.6469
$ws.Range("H137").Value = 3115.99
$ws.Range("I137").Value = 1250.2632
$ws.Range("K137").Value = 3750.7896
$ws.Range("M137").Value = -1200.7896
$ws.Range("H141").Value = 1790.2
$ws.Range("I141").Value = 1790.2
$ws.Range("K141").Value = 5370.6
$ws.Range("M141").Value = -190.6000000000004
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 8293.111000000001
$ws.Range("I74").Value = 6888.6665
$ws.Range("K74").Value = 6888.6665
$ws.Range("M74").Value = -6014.6665
$ws.Range("H77").Value = 8293.111000000001
$ws.Range("I77").Value = 6888.6665
$ws.Range("K77").Value = 34443.3325
$ws.Range("M77").Value = -30075.3325
$ws.Range("H102").Value = 9357.296
$ws.Range("I102").Value = 6244.7144
$ws.Range("K102").Value = 6244.7144
$ws.Range("M102").Value = -4622.7144
$ws.Range("H110").Value = 41674476
$ws.Range("I110").Value = 47620260
$ws.Range("K110").Value = 47620260
$ws.Range("M110").Value = -47618215
$ws.Range("H122").Value = 9525878
$ws.Range("I122").Value = 11495835
$ws.Range("J122").Value = 4416.6665
$ws.Range("K122").Value = 34487505
$ws.Range("L122").Value = 13249.9995
$ws.Range("M122").Value = -34485055
$ws.Range("N122").Value = -18149.9995
$ws.Range("H132").Value = 200003800
$ws.Range("J132").Value = 4534.5
$ws.Range("L132").Value = 13603.5
$ws.Range("N132").Value = -18663.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 61659
$ws.Range("J35").Value = 61659
$ws.Range("L35").Value = 61659
$ws.Range("N35").Value = -62279
$ws.Range("H86").Value = 23811848
$ws.Range("I86").Value = 41668540
$ws.Range("J86").Value = 2927.4443
$ws.Range("K86").Value = 41668540
$ws.Range("L86").Value = 2927.4443
$ws.Range("M86").Value = -41667417
$ws.Range("N86").Value = -5173.4443
$ws.Range("H89").Value = 23811848
$ws.Range("I89").Value = 41668540
$ws.Range("J89").Value = 2927.4443
$ws.Range("K89").Value = 208342700
$ws.Range("L89").Value = 14637.2215
$ws.Range("M89").Value = -208337084
$ws.Range("N89").Value = -25869.2215
$ws.Range("H107").Value = 16687472
$ws.Range("I107").Value = 9147.087
$ws.Range("J107").Value = 71487680
$ws.Range("K107").Value = 9147.087
$ws.Range("L107").Value = 71487680
$ws.Range("M107").Value = -7227.087
$ws.Range("N107").Value = -71491520
$ws.Range("H122").Value = 49999
$ws.Range("J122").Value = 49999
$ws.Range("L122").Value = 49999
$ws.Range("N122").Value = -59799
$ws.Range("H134").Value = 2889.7856
$ws.Range("I134").Value = 2627.875
$ws.Range("J134").Value = 4461.25
$ws.Range("K134").Value = 7883.625
$ws.Range("L134").Value = 13383.75
$ws.Range("M134").Value = -5348.625
$ws.Range("N134").Value = -18453.75
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 90.25
$ws.Range("I7").Value = 88.40000000000001
$ws.Range("J7").Value = 99.5
$ws.Range("K7").Value = 88.40000000000001
$ws.Range("L7").Value = 99.5
$ws.Range("M7").Value = 24.59999999999999
$ws.Range("N7").Value = -325.5
$ws.Range("H31").Value = 3855.3635
$ws.Range("I31").Value = 1193.625
$ws.Range("K31").Value = 1193.625
$ws.Range("M31").Value = -898.625
$ws.Range("H34").Value = 3855.3635
$ws.Range("I34").Value = 1193.625
$ws.Range("K34").Value = 1193.625
$ws.Range("M34").Value = -991.625
$ws.Range("H99").Value = 2390.5
$ws.Range("I99").Value = 2106.1667
$ws.Range("K99").Value = 2106.1667
$ws.Range("M99").Value = -608.1667000000002
$ws.Range("H105").Value = 1873.6
$ws.Range("I105").Value = 1869.75
$ws.Range("J105").Value = 1889
$ws.Range("K105").Value = 1869.75
$ws.Range("L105").Value = 1889
$ws.Range("M105").Value = -122.75
$ws.Range("N105").Value = -5383
$ws.Range("H126").Value = 2390.5
$ws.Range("I126").Value = 2106.1667
$ws.Range("K126").Value = 6318.500100000001
$ws.Range("M126").Value = -3848.500100000001
$ws.Range("H132").Value = 1431490.4
$ws.Range("I132").Value = 1002586.6
$ws.Range("K132").Value = 3007759.8
$ws.Range("M132").Value = -3005229.8
$ws.Range("H134").Value = 3933.111
$ws.Range("I134").Value = 3000
$ws.Range("J134").Value = 4199.7144
$ws.Range("K134").Value = 9000
$ws.Range("L134").Value = 12599.1432
$ws.Range("M134").Value = -6465
$ws.Range("N134").Value = -17669.1432
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1745
$ws.Range("I3").Value = 1745
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 5235
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = $null
$ws.Range("N3").Value = -5123
$ws.Range("H113").Value = 747.7857
$ws.Range("I113").Value = 326.2
$ws.Range("K113").Value = 978.5999999999999
$ws.Range("M113").Value = 1191.4
$ws.Range("H138").Value = 3820155.5
$ws.Range("I138").Value = 6001049
$ws.Range("J138").Value = 185333
$ws.Range("K138").Value = 18003147
$ws.Range("L138").Value = 555999
$ws.Range("M138").Value = -17998007
$ws.Range("N138").Value = -566279
$ws.Range("H139").Value = 1236274.4
$ws.Range("I139").Value = 1755695.1
$ws.Range("J139").Value = 2650
$ws.Range("K139").Value = 5267085.300000001
$ws.Range("L139").Value = 7950
$ws.Range("M139").Value = -5261945.300000001
$ws.Range("N139").Value = -18230
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 6000
$ws.Range("I48").Value = 6000
$ws.Range("K48").Value = 6000
$ws.Range("M48").Value = -5515
$ws.Range("H80").Value = 3650.25
$ws.Range("I80").Value = 3761.8
$ws.Range("J80").Value = 3599.5454
$ws.Range("K80").Value = 3761.8
$ws.Range("L80").Value = 3599.5454
$ws.Range("M80").Value = -2763.8
$ws.Range("N80").Value = -5595.5454
$ws.Range("H83").Value = 3650.25
$ws.Range("I83").Value = 3761.8
$ws.Range("J83").Value = 3599.5454
$ws.Range("K83").Value = 18809
$ws.Range("L83").Value = 17997.727
$ws.Range("M83").Value = -13817
$ws.Range("N83").Value = -27981.727
$ws.Range("H102").Value = 1308.3334
$ws.Range("I102").Value = 1308.3334
$ws.Range("K102").Value = 1308.3334
$ws.Range("M102").Value = 313.6666
$ws.Range("H122").Value = 33335670
$ws.Range("I122").Value = 2314.6155
$ws.Range("K122").Value = 6943.8465
$ws.Range("M122").Value = -4493.8465
$ws.Range("H123").Value = 43560
$ws.Range("J123").Value = 43560
$ws.Range("L123").Value = 43560
$ws.Range("N123").Value = -48460
$ws.Range("H132").Value = 203630.4
$ws.Range("I132").Value = 297478
$ws.Range("J132").Value = 4204.25
$ws.Range("K132").Value = 892434
$ws.Range("L132").Value = 12612.75
$ws.Range("M132").Value = -889904
$ws.Range("N132").Value = -17672.75
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4145
$ws.Range("I16").Value = 4025
$ws.Range("K16").Value = 4025
$ws.Range("M16").Value = -3855
$ws.Range("H22").Value = 5380343.5
$ws.Range("I22").Value = 4666.3335
$ws.Range("J22").Value = 10756021
$ws.Range("K22").Value = 4666.3335
$ws.Range("L22").Value = 10756021
$ws.Range("M22").Value = -4371.3335
$ws.Range("N22").Value = -10756611
$ws.Range("H27").Value = 5380343.5
$ws.Range("I27").Value = 4666.3335
$ws.Range("J27").Value = 10756021
$ws.Range("K27").Value = 4666.3335
$ws.Range("L27").Value = 10756021
$ws.Range("M27").Value = -4559.3335
$ws.Range("N27").Value = -10756235
$ws.Range("H35").Value = 3466
$ws.Range("I35").Value = 1943.3334
$ws.Range("J35").Value = 5750
$ws.Range("K35").Value = 1943.3334
$ws.Range("L35").Value = 5750
$ws.Range("M35").Value = -1607.3334
$ws.Range("N35").Value = -6422
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").Value = $null
$ws.Range("H46").Value = 2962.1738
$ws.Range("I46").Value = 2432.5
$ws.Range("J46").Value = 3149.1177
$ws.Range("K46").Value = 2432.5
$ws.Range("L46").Value = 3149.1177
$ws.Range("M46").Value = -2244.5
$ws.Range("N46").Value = -3525.1177
$ws.Range("H122").Value = 3037.1853
$ws.Range("I122").Value = 2521.913
$ws.Range("J122").Value = 6000
$ws.Range("K122").Value = 7565.739
$ws.Range("L122").Value = 18000
$ws.Range("M122").Value = -5115.739
$ws.Range("N122").Value = -22900
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 20000
$ws.Range("I37").Value = 20000
$ws.Range("K37").Value = 20000
$ws.Range("M37").Value = -19797
$ws.Range("H136").Value = 3986.4092
$ws.Range("I136").Value = 2821.1
$ws.Range("K136").Value = 8463.299999999999
$ws.Range("M136").Value = -5913.299999999999
